$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update F column ("想去人数") values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 12
$ws1.Range("F7").Value = 1599
$ws1.Range("F9").Value = 14
$ws1.Range("F10").Value = 1342
$ws1.Range("F11").Value = 119
$ws1.Range("F12").Value = 14
$ws1.Range("F13").Value = 233
$ws1.Range("F14").Value = 175
$ws1.Range("F18").Value = 243
$ws1.Range("F20").Value = 202

# Sheet "全部类型" (fourth sheet) - update F column ("想去人数") values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 12
$ws4.Range("F7").Value = 1599
$ws4.Range("F10").Value = 14
$ws4.Range("F11").Value = 1342
$ws4.Range("F12").Value = 119
$ws4.Range("F13").Value = 15
$ws4.Range("F14").Value = 233
$ws4.Range("F15").Value = 175
$ws4.Range("F19").Value = 243
$ws4.Range("F21").Value = 202
